$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Cell B11 currently holds the text "R40" (row 11 is the "R40 .. Good
# Night" rule). The edit replaces its content with the literal text "1"
# (a brand-new shared string), keeping the cell's existing formatting
# untouched.
#
# Assigning a numeric-looking string straight to .Value would make Excel
# store it as the *number* 1 instead of the text "1", and prefixing it
# with a leading apostrophe (the usual way to force text) would stamp the
# cell with a "quote prefix" style, changing its style id. So: stash the
# current formatting in a scratch cell, write the text value, then paste
# the stashed formatting back on top.
$scratch = $ws.Range("Z1000")

$ws.Range("B11").Copy()
$scratch.PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("B11").Value = "'1"  # leading apostrophe forces literal text "1"

$scratch.Copy()
$ws.Range("B11").PasteSpecial(-4122)   # xlPasteFormats

$scratch.Clear()
$excel.CutCopyMode = $false
